# regional_fuel_type_defaults.xlsx
# "create separate columns for system 4 and 6"
#
# Column J ("heating_coil_type_sys4and6") is split into two columns:
#   J = heating_coil_type_sys4 (keeps the existing values)
#   K = heating_coil_type_sys6 (new column, duplicate of the sys4 values)
# The old column K ("fan_type" / var_speed_drive) shifts right to column L.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new blank column at K, shifting fan_type (old K) to L.
$ws.Columns.Item(11).Insert()

# Split the header: J keeps "sys4", new K gets "sys6".
$ws.Range("J1").Value = "heating_coil_type_sys4"
$ws.Range("K1").Value = "heating_coil_type_sys6"

# Duplicate the sys4 data (values + formatting) into the new sys6 column.
$ws.Range("J2:J81").Copy()
$ws.Range("K2:K81").PasteSpecial(-4104)

# Match the new column's width to its sys4 neighbour.
$ws.Columns.Item(11).ColumnWidth = $ws.Columns.Item(10).ColumnWidth

$excel.CutCopyMode = $false

[void]$ws.Range("K54").Select()
